# Auto-generated Excel COM-interop script
# Applies numeric cell updates across multiple sheets per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 50000
$ws.Range("J47").Value = 50000
$ws.Range("L47").Value = 50000
$ws.Range("N47").Value = -51944

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1909
$ws.Range("J112").Value = 1420
$ws.Range("L112").Value = 4260
$ws.Range("N112").Value = -6476

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3923489.2
$ws.Range("I132").Value = 4349667.5
$ws.Range("J132").Value = 2651.2
$ws.Range("K132").Value = 13049002.5
$ws.Range("L132").Value = 7953.599999999999
$ws.Range("M132").Value = -13046472.5
$ws.Range("N132").Value = -13013.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3610.1045
$ws.Range("I138").Value = 1330.4419
$ws.Range("J138").Value = 7694.5
$ws.Range("K138").Value = 3991.3257
$ws.Range("L138").Value = 23083.5
$ws.Range("M138").Value = 1148.6743
$ws.Range("N138").Value = -33363.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6420.58
$ws.Range("I32").Value = 6341.0425
$ws.Range("J32").Value = 7666.6665
$ws.Range("K32").Value = 6341.0425
$ws.Range("L32").Value = 7666.6665
$ws.Range("M32").Value = -6054.0425
$ws.Range("N32").Value = -8240.666499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1728.1082
$ws.Range("I122").Value = 1483.92
$ws.Range("J122").Value = 2236.8333
$ws.Range("K122").Value = 4451.76
$ws.Range("L122").Value = 6710.499899999999
$ws.Range("M122").Value = -2001.76
$ws.Range("N122").Value = -11610.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1713.2142
$ws.Range("I99").Value = 1416.1305
$ws.Range("K99").Value = 1416.1305
$ws.Range("M99").Value = 81.86950000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2953.5366
$ws.Range("I134").Value = 3207.9119
$ws.Range("J134").Value = 1718
$ws.Range("K134").Value = 9623.735700000001
$ws.Range("L134").Value = 5154
$ws.Range("M134").Value = -7088.735700000001
$ws.Range("N134").Value = -10224

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1722.4166
$ws.Range("I16").Value = 1721.25
$ws.Range("J16").Value = 1724.75
$ws.Range("K16").Value = 1721.25
$ws.Range("L16").Value = 1724.75
$ws.Range("M16").Value = -1434.25
$ws.Range("N16").Value = -2298.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12503171
$ws.Range("I58").Value = 2210.6667
$ws.Range("K58").Value = 2210.6667
$ws.Range("M58").Value = -2007.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1722.4166
$ws.Range("I113").Value = 1721.25
$ws.Range("J113").Value = 1724.75
$ws.Range("K113").Value = 1721.25
$ws.Range("L113").Value = 1724.75
$ws.Range("M113").Value = 448.75
$ws.Range("N113").Value = -6064.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H114").Value = 31000
$ws.Range("J114").Value = 31000
$ws.Range("L114").Value = 31000
$ws.Range("N114").Value = -39678

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2338.2173
$ws.Range("I132").Value = 3336.5454
$ws.Range("J132").Value = 1423.0834
$ws.Range("K132").Value = 10009.6362
$ws.Range("L132").Value = 4269.2502
$ws.Range("M132").Value = -7479.636200000001
$ws.Range("N132").Value = -9329.2502

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3068.6428
$ws.Range("I134").Value = 3362.4666
$ws.Range("J134").Value = 2334.0833
$ws.Range("K134").Value = 10087.3998
$ws.Range("L134").Value = 7002.249899999999
$ws.Range("M134").Value = -7552.399800000001
$ws.Range("N134").Value = -12072.2499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 12503171
$ws.Range("I136").Value = 2210.6667
$ws.Range("K136").Value = 6632.000100000001
$ws.Range("M136").Value = -4082.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 411.3143
$ws.Range("I12").Value = 28.8
$ws.Range("J12").Value = 564.3200000000001
$ws.Range("K12").Value = 86.40000000000001
$ws.Range("L12").Value = 1692.96
$ws.Range("M12").Value = 86.59999999999999
$ws.Range("N12").Value = -2038.96

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 6140.5
$ws.Range("I34").Value = 250
$ws.Range("J34").Value = 10960
$ws.Range("K34").Value = 750
$ws.Range("L34").Value = 32880
$ws.Range("M34").Value = -666
$ws.Range("N34").Value = -33048

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 750.6875
$ws.Range("J114").Value = 1966.2
$ws.Range("L114").Value = 5898.6
$ws.Range("N114").Value = -12406.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1613.1
$ws.Range("I117").Value = 649.75
$ws.Range("J117").Value = 2255.3333
$ws.Range("K117").Value = 1949.25
$ws.Range("L117").Value = 6765.999899999999
$ws.Range("M117").Value = 1492.75
$ws.Range("N117").Value = -13649.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 27430
$ws.Range("I129").Value = 3114.4443
$ws.Range("J129").Value = 47324.547
$ws.Range("K129").Value = 9343.332900000001
$ws.Range("L129").Value = 141973.641
$ws.Range("M129").Value = -4343.332900000001
$ws.Range("N129").Value = -151973.641

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1231.375
$ws.Range("I131").Value = 1917.5
$ws.Range("J131").Value = 1039.26
$ws.Range("K131").Value = 5752.5
$ws.Range("L131").Value = 3117.78
$ws.Range("M131").Value = -712.5
$ws.Range("N131").Value = -13197.78

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2131.25
$ws.Range("I61").Value = 563.8889
$ws.Range("K61").Value = 563.8889
$ws.Range("M61").Value = -361.8889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2131.25
$ws.Range("I113").Value = 563.8889
$ws.Range("K113").Value = 563.8889
$ws.Range("M113").Value = 1606.1111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3929.5881
$ws.Range("J136").Value = 3833.3333
$ws.Range("L136").Value = 11499.9999
$ws.Range("N136").Value = -16599.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 39999
$ws.Range("J115").Value = 39999
$ws.Range("L115").Value = 39999
$ws.Range("N115").Value = -43133

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5279.973
$ws.Range("I132").Value = 2238.5334
$ws.Range("J132").Value = 18314.715
$ws.Range("K132").Value = 6715.600199999999
$ws.Range("L132").Value = 54944.145
$ws.Range("M132").Value = -4185.600199999999
$ws.Range("N132").Value = -60004.145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3342.353
$ws.Range("I136").Value = 3474.818
$ws.Range("K136").Value = 10424.454
$ws.Range("M136").Value = -7874.454000000002
